$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IAM_PERMISSION")

$ws.Range("G8").Value = "/manager/site-statistics"
$ws.Range("G9").Value = "/manager/microservice"
$ws.Range("G10").Value = "/manager/instance"
$ws.Range("G11").Value = "/manager/configuration"
$ws.Range("G12").Value = "/manager/route"
$ws.Range("G13").Value = "/manager/api-test"
$ws.Range("G14").Value = "/manager/api-test"
$ws.Range("G15").Value = "/manager/api-overview"
$ws.Range("G16").Value = "/manager/api-overview"
$ws.Range("G17").Value = "/manager/site-statistics"
